# Refitting NCDEs to individual patients (for manuscript figure)
# Adds a "Label" column (H) indicating Control (0) vs MDD (1) patient group,
# and refreshes the refit metrics in columns D/E/F that changed slightly
# with the new fit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for column H, styled like the other header cells (B1:G1)
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Label"

# Refreshed metric values (columns D/E/F) from the re-fit run
$ws.Range("D2").Value = 0.5924117134068744
$ws.Range("E2").Value = 0.5924117134068744

$ws.Range("D3").Value = 0.9997617181339506
$ws.Range("E3").Value = 0.9997617181339506

$ws.Range("D4").Value = 0.3885816417400595
$ws.Range("E4").Value = 0.3885816417400595

$ws.Range("D5").Value = 0.908811407480433
$ws.Range("E5").Value = 0.908811407480433

$ws.Range("D6").Value = 0.4707147033387259
$ws.Range("E6").Value = 0.4707147033387259

$ws.Range("D7").Value = 0.5785510606564803
$ws.Range("E7").Value = 0.4214489393435197

$ws.Range("D8").Value = 0.4565525211171457
$ws.Range("E8").Value = 0.5434474788828543

$ws.Range("D9").Value = 0.4269486203750447
$ws.Range("E9").Value = 0.5730513796249552

$ws.Range("D10").Value = 0.5025401433020892
$ws.Range("E10").Value = 0.4974598566979108

$ws.Range("D11").Value = 0.9861414006191495
$ws.Range("E11").Value = 0.01385859938085054
$ws.Range("F11").Value = 1.564698815345764

# New "Label" column values: 0 = Control patient, 1 = MDD patient
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("H11").Value = 1
$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 1
$ws.Range("H21").Value = 1
